# Apply the edits described by the commit:
#   "added power analysis and fixed things up before submission"
#
# 1. Label the best (lowest-AIC) model in the "prey" question block (row 9,
#    repeated in the summary block at row 16) and the best model in the
#    "SST" question block (row 15) with a bold "(best ... model)" suffix
#    appended to the existing label text.
# 2. Re-point the stray $L$20 / $L$26 anchors in the column-E "ΔAIC"
#    formulas onto the actual lowest-AIC row of each block ($D$4, $D$9,
#    $D$15) and add formulas to the rows that previously held only a
#    literal value.
# 3. Leave the window scrolled/selected roughly where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Append bold "(best ... model)" suffixes -------------------------

function Set-LabelWithBoldSuffix($cell, [string]$base, [string]$suffix) {
    $full = $base + $suffix
    $cell.Value = $full
    $text = $cell.Text
    $start = $text.Length - $suffix.Length + 1
    $cell.Characters($start, $suffix.Length).Font.Bold = $true
}

Set-LabelWithBoldSuffix $ws.Range("A9")  "Primary prey (fixed effects) " "(best prey model)"
Set-LabelWithBoldSuffix $ws.Range("A16") "Primary prey (fixed effects) " "(best prey model)"
Set-LabelWithBoldSuffix $ws.Range("A15") "SST (fixed effects) "          "(best SST model)"

# --- 2. Fix up the ΔAIC formulas in column E -----------------------------

# Block 1 (rows 4-7): anchor on row 4, the lowest-AIC row of the block.
$ws.Range("E4").Formula = "=D4-`$D`$4"
$ws.Range("E5:E7").Formula = "=D5-`$D`$4"

# Block 2 (rows 9-13): anchor on row 9, the lowest-AIC row of the block.
$ws.Range("E9").Formula = "=D9-`$D`$9"
$ws.Range("E10:E13").Formula = "=D10-`$D`$9"

# Block 3 (rows 15-21): anchor on row 15, the lowest-AIC row of the block.
$ws.Range("E15").Formula = "=D15-`$D`$15"
$ws.Range("E16").Formula = "=D16-`$D`$15"
$ws.Range("E17").Formula = "=D17-`$D`$15"
$ws.Range("E18").Formula = "=D18-`$D`$15"
$ws.Range("E19").Formula = "=D19-`$D`$15"
$ws.Range("E20").Formula = "=D20-`$D`$15"
$ws.Range("E21").Formula = "=D21-`$D`$15"

# --- 3. Window scroll position / selection -------------------------------

$win = $ws.Application.ActiveWindow
$win.ScrollRow = 8
$ws.Range("E21").Select()

$wb.Save()
